$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.517.70'
$ws.Range('E2').Value = '  +2.07%  '
$ws.Range('D3').Value = '1.911.22'
$ws.Range('E3').Value = '  +5.45%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5055'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.69%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3957'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09821'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.163'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.66'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.555'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.15'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.23%  '
$ws.Range('D14').Value = '1.920.47'
$ws.Range('E14').Value = '  +5.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.582'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9998'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.91'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06652'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('E20').Value = '  +5.68%  '
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.298'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.44%  '
$ws.Range('D23').Value = '28.565.03'
$ws.Range('E23').Value = '  +2.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.47'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.281'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.741'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +14.33%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.394'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.95%  '
$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D28').Value = '2.127.12'
$ws.Range('E28').Value = '  +5.29%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '21.38'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.03%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '159.44'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '129.47'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.07%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.107'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.75%  '
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.1072'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.715'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.73%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.636'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.904'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.86%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06803'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.46%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02448'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.17%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.276'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +10.10%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2233'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.58%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.104'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6458'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.42%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.66%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.193'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.20%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9995'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.15%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.69'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.63%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6114'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.04%  '
$ws.Range('B48').Value = 'WEMIXTOKEN'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.283'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('B49').Value = 'PancakeSwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.662'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.055'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.15%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '125.40'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.94%  '
